$wb = $excel.ActiveWorkbook

# --- Sheet "all_tools": update num_warnings (G) and recomputed correlation stats (I-L) for dataset rows 9_gc/9_bc/9_nc (rows 13-24) ---
$wsAllTools = $wb.Worksheets.Item("all_tools")
$wsAllTools.Range("G13").Value = 73
$wsAllTools.Range("I13").Value = -0.4535394202249742
$wsAllTools.Range("J13").Value = 0.08070214265077501
$wsAllTools.Range("K13").Value = -0.5512130501182249
$wsAllTools.Range("L13").Value = 0.09862398123473154

$wsAllTools.Range("G14").Value = 73
$wsAllTools.Range("I14").Value = -0.02387049580131443
$wsAllTools.Range("J14").Value = 0.92675547372309
$wsAllTools.Range("K14").Value = -0.09290107586262218
$wsAllTools.Range("L14").Value = 0.7985237548304135

$wsAllTools.Range("G15").Value = 73
$wsAllTools.Range("I15").Value = 0.09656090991705352
$wsAllTools.Range("J15").Value = 0.7120793980044939
$wsAllTools.Range("K15").Value = 0.1428819499477476
$wsAllTools.Range("L15").Value = 0.6937488280957302

$wsAllTools.Range("G16").Value = 73
$wsAllTools.Range("I16").Value = -0.7399853698407473
$wsAllTools.Range("J16").Value = 0.004375235749920733
$wsAllTools.Range("K16").Value = -0.8670767080511405
$wsAllTools.Range("L16").Value = 0.001159768265192443

$wsAllTools.Range("G17").Value = 73
$wsAllTools.Range("I17").Value = -0.3103164454170876
$wsAllTools.Range("J17").Value = 0.2320634889020341
$wsAllTools.Range("K17").Value = -0.4149581388530457
$wsAllTools.Range("L17").Value = 0.233093730241445

$wsAllTools.Range("G18").Value = 73
$wsAllTools.Range("I18").Value = 0.4535394202249742
$wsAllTools.Range("J18").Value = 0.08070214265077501
$wsAllTools.Range("K18").Value = 0.5635998602332413
$wsAllTools.Range("L18").Value = 0.08974939558201678

$wsAllTools.Range("G19").Value = 73
$wsAllTools.Range("I19").Value = 0.4828045495852675
$wsAllTools.Range("J19").Value = 0.06499039472076076
$wsAllTools.Range("K19").Value = 0.5932707052178214
$wsAllTools.Range("L19").Value = 0.0706202179291523

$wsAllTools.Range("G20").Value = 73
$wsAllTools.Range("I20").Value = -0.3580574370197164
$wsAllTools.Range("J20").Value = 0.1679207532945924
$wsAllTools.Range("K20").Value = -0.4954724046006516
$wsAllTools.Range("L20").Value = 0.1453294522910624

$wsAllTools.Range("G21").Value = 73
$wsAllTools.Range("I21").Value = -0.4535394202249742
$wsAllTools.Range("J21").Value = 0.08070214265077501
$wsAllTools.Range("K21").Value = -0.6007602905782901
$wsAllTools.Range("L21").Value = 0.06625295074379814

$wsAllTools.Range("G22").Value = 73
$wsAllTools.Range("I22").Value = 0.1193524790065721
$wsAllTools.Range("J22").Value = 0.6457756768370824
$wsAllTools.Range("K22").Value = 0.1486417213801955
$wsAllTools.Range("L22").Value = 0.6819355638686473

$wsAllTools.Range("G23").Value = 73
$wsAllTools.Range("I23").Value = 0.167093470609201
$wsAllTools.Range("J23").Value = 0.5199036173455835
$wsAllTools.Range("K23").Value = 0.2725098225303584
$wsAllTools.Range("L23").Value = 0.446215643690079

$wsAllTools.Range("G24").Value = 73
$wsAllTools.Range("I24").Value = -0.3580574370197164
$wsAllTools.Range("J24").Value = 0.1679207532945924
$wsAllTools.Range("K24").Value = -0.4706987843706191
$wsAllTools.Range("L24").Value = 0.1697475039817557

# --- Sheet "typestate_checker": update num_snippets_warnings (F), num_warnings (G), and newly computed correlation stats (I-L) for dataset rows 9_gc/9_bc/9_nc (rows 13-24) ---
$wsTypestate = $wb.Worksheets.Item("typestate_checker")
$wsTypestate.Range("F13").Value = 10
$wsTypestate.Range("G13").Value = 37
$wsTypestate.Range("I13").Value = -0.4787549991450212
$wsTypestate.Range("J13").Value = 0.07217560549492458
$wsTypestate.Range("K13").Value = -0.6292853089020909
$wsTypestate.Range("L13").Value = 0.05124855216842294

$wsTypestate.Range("F14").Value = 10
$wsTypestate.Range("G14").Value = 37
$wsTypestate.Range("I14").Value = -0.2267786838055363
$wsTypestate.Range("J14").Value = 0.3943870594034554
$wsTypestate.Range("K14").Value = -0.2860387767736777
$wsTypestate.Range("L14").Value = 0.4230203924441358

$wsTypestate.Range("F15").Value = 10
$wsTypestate.Range("G15").Value = 37
$wsTypestate.Range("I15").Value = -0.1019294382875251
$wsTypestate.Range("J15").Value = 0.7040542681897126
$wsTypestate.Range("K15").Value = -0.0765092055676006
$wsTypestate.Range("L15").Value = 0.8336123677972922

$wsTypestate.Range("F16").Value = 10
$wsTypestate.Range("G16").Value = 37
$wsTypestate.Range("I16").Value = -0.579545525280815
$wsTypestate.Range("J16").Value = 0.02951512807757192
$wsTypestate.Range("K16").Value = -0.7119187333033755
$wsTypestate.Range("L16").Value = 0.02091481468718881

$wsTypestate.Range("F17").Value = 10
$wsTypestate.Range("G17").Value = 37
$wsTypestate.Range("I17").Value = -0.1259881576697424
$wsTypestate.Range("J17").Value = 0.6360988735986226
$wsTypestate.Range("K17").Value = -0.1906925178491184
$wsTypestate.Range("L17").Value = 0.5977007516614028

$wsTypestate.Range("F18").Value = 10
$wsTypestate.Range("G18").Value = 37
$wsTypestate.Range("I18").Value = 0.579545525280815
$wsTypestate.Range("J18").Value = 0.02951512807757192
$wsTypestate.Range("K18").Value = 0.7437008196115621
$wsTypestate.Range("L18").Value = 0.01366958411527145

$wsTypestate.Range("F19").Value = 10
$wsTypestate.Range("G19").Value = 37
$wsTypestate.Range("I19").Value = 0.560611910581388
$wsTypestate.Range("J19").Value = 0.0366903087793031
$wsTypestate.Range("K19").Value = 0.7172738021962557
$wsTypestate.Range("L19").Value = 0.01954204435368506

$wsTypestate.Range("F20").Value = 10
$wsTypestate.Range("G20").Value = 37
$wsTypestate.Range("I20").Value = -0.4787549991450212
$wsTypestate.Range("J20").Value = 0.07217560549492458
$wsTypestate.Range("K20").Value = -0.6419981434253655
$wsTypestate.Range("L20").Value = 0.04536158917864154

$wsTypestate.Range("F21").Value = 10
$wsTypestate.Range("G21").Value = 37
$wsTypestate.Range("I21").Value = -0.2771739468734333
$wsTypestate.Range("J21").Value = 0.2978975979923409
$wsTypestate.Range("K21").Value = -0.3750286184365996
$wsTypestate.Range("L21").Value = 0.2855969029688312

$wsTypestate.Range("F22").Value = 10
$wsTypestate.Range("G22").Value = 37
$wsTypestate.Range("I22").Value = 0.1259881576697424
$wsTypestate.Range("J22").Value = 0.6360988735986226
$wsTypestate.Range("K22").Value = 0.1461975970176575
$wsTypestate.Range("L22").Value = 0.6869410188538527

$wsTypestate.Range("F23").Value = 10
$wsTypestate.Range("G23").Value = 37
$wsTypestate.Range("I23").Value = 0.1259881576697424
$wsTypestate.Range("J23").Value = 0.6360988735986226
$wsTypestate.Range("K23").Value = 0.2161181868956676
$wsTypestate.Range("L23").Value = 0.5487107060733141

$wsTypestate.Range("F24").Value = 10
$wsTypestate.Range("G24").Value = 37
$wsTypestate.Range("I24").Value = -0.2267786838055363
$wsTypestate.Range("J24").Value = 0.3943870594034554
$wsTypestate.Range("K24").Value = -0.3114644458202268
$wsTypestate.Range("L24").Value = 0.3810089567050594

